$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.200.50"
$ws.Range("E2").Value = '  +0.07%  '

$ws.Range("D3").Value = "'1.851.85"
$ws.Range("E3").Value = '  -0.25%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = "'0.6976"
$ws.Range("E5").Value = '  +1.17%  '

$ws.Range("D6").Value = "'237.60"
$ws.Range("E6").Value = '  -0.21%  '

$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = '  +0.11%  '

$ws.Range("D8").Value = "'0.07882"
$ws.Range("E8").Value = '  +2.41%  '

$ws.Range("D9").Value = "'0.3014"
$ws.Range("E9").Value = '  -1.26%  '

$ws.Range("D10").Value = "'23.57"
$ws.Range("E10").Value = '  +1.38%  '

$ws.Range("D11").Value = "'0.08117"
$ws.Range("E11").Value = '  +0.65%  '

$ws.Range("D12").Value = "'1.846.55"
$ws.Range("E12").Value = '  -1.08%  '

$ws.Range("D13").Value = "'5.186"
$ws.Range("E13").Value = '  -0.13%  '

$ws.Range("D14").Value = "'0.7042"
$ws.Range("E14").Value = '  -2.43%  '

$ws.Range("D15").Value = "'89.47"
$ws.Range("E15").Value = '  +0.13%  '

$ws.Range("D16").Value = "'29.242.05"
$ws.Range("E16").Value = '  +0.21%  '

$ws.Range("D17").Value = "'5.803"
$ws.Range("E17").Value = '  +1.18%  '

$ws.Range("D18").Value = "'0.000007800"
$ws.Range("E18").Value = '  +0.12%  '

$ws.Range("D19").Value = "'13.18"
$ws.Range("E19").Value = '  -0.50%  '

$ws.Range("D20").Value = "'235.55"
$ws.Range("E20").Value = '  +0.23%  '

$ws.Range("E21").Value = '  +0.09%  '

$ws.Range("D22").Value = "'2.105.84"
$ws.Range("E22").Value = '  -0.26%  '

$ws.Range("E23").Value = '  -0.04%  '

$ws.Range("D24").Value = "'7.520"
$ws.Range("E24").Value = '  +1.00%  '

$ws.Range("D25").Value = "'162.65"
$ws.Range("E25").Value = '  +0.47%  '

$ws.Range("D26").Value = "'8.847"
$ws.Range("E26").Value = '  -1.33%  '

$ws.Range("D27").Value = "'0.1412"
$ws.Range("E27").Value = '  -1.28%  '

$ws.Range("E28").Value = '  -0.30%  '

$ws.Range("D29").Value = "'1.910"
$ws.Range("E29").Value = '  -2.22%  '

$ws.Range("D30").Value = "'1.408"
$ws.Range("E30").Value = '  +0.62%  '

$ws.Range("D31").Value = "'1.471"
$ws.Range("E31").Value = '  -1.03%  '

$ws.Range("D32").Value = "'4.320"
$ws.Range("E32").Value = '  -4.70%  '

$ws.Range("E33").Value = '  +0.04%  '

$ws.Range("D34").Value = "'0.05145"
$ws.Range("E34").Value = '  -0.87%  '

$ws.Range("E35").Value = '  -1.96%  '

$ws.Range("D36").Value = "'0.7101"
$ws.Range("E36").Value = '  +0.77%  '

$ws.Range("D37").Value = "'0.9991"
$ws.Range("E37").Value = '  -2.08%  '

$ws.Range("E38").Value = '  +0.24%  '

$ws.Range("D39").Value = "'0.01841"
$ws.Range("E39").Value = '  -0.44%  '

$ws.Range("D40").Value = "'2.707"
$ws.Range("E40").Value = '  +0.96%  '

$ws.Range("D41").Value = "'1.153.58"
$ws.Range("E41").Value = '  +5.32%  '

$ws.Range("E42").Value = '  -0.23%  '

$ws.Range("D43").Value = "'5.973"
$ws.Range("E43").Value = '  -0.24%  '

$ws.Range("D44").Value = "'0.4234"
$ws.Range("E44").Value = '  -1.18%  '

$ws.Range("D45").Value = "'69.92"
$ws.Range("E45").Value = '  -0.93%  '

$ws.Range("E46").Value = '  -0.01%  '

$ws.Range("D47").Value = "'102.68"
$ws.Range("E47").Value = '  +0.41%  '

$ws.Range("D48").Value = "'0.5298"
$ws.Range("E48").Value = '  -2.89%  '

$ws.Range("D49").Value = "'1.734"
$ws.Range("E49").Value = '  -3.03%  '

$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = "'1.993.11"
$ws.Range("E50").Value = '  -0.55%  '

$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = "'9.117"
$ws.Range("E51").Value = '  -0.52%  '
